$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "If using Project Scarlett, set the active solution platform to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If using an Xbox Series X|S devkit, set the active solution platform to ",
    2
)
